$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92, shifting the existing data (old rows
# 92..184) down to 93..185, and fill the new row with this week's record.
$ws.Rows.Item(92).Insert()

$ws.Cells.Item(92, 1).Value  = 7
$ws.Cells.Item(92, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(92, 3).Value  = "Ñuble"
$ws.Cells.Item(92, 4).Value  = 45167
$ws.Cells.Item(92, 5).Value  = 16
$ws.Cells.Item(92, 6).Value  = "Fruta"
$ws.Cells.Item(92, 7).Value  = 100108
$ws.Cells.Item(92, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(92, 9).Value  = 100108002
$ws.Cells.Item(92, 10).Value = "Mango"
$ws.Cells.Item(92, 11).Value = "Sin especificar"
$ws.Cells.Item(92, 12).Value = "Primera"
$ws.Cells.Item(92, 13).Value = 50
$ws.Cells.Item(92, 14).Value = 9000
$ws.Cells.Item(92, 15).Value = 9000
$ws.Cells.Item(92, 16).Value = 9000
$ws.Cells.Item(92, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(92, 18).Value = "Brasil"
$ws.Cells.Item(92, 19).Value = 2250
$ws.Cells.Item(92, 20).Value = 4
